# Auto-generated: add rows 6-13 (Waktu/Input/Hasil) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bytes = [Convert]::FromBase64String("MjAyNS0wOC0wMyAyMDozMjoyNA==")
$ws.Range("A6").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("bGlzdCBzZWtvbGFoIHlhbmcgaWt1dCBkaSBhcGxpa2FzaSBzYWF0IGluaQ==")
$ws.Range("B6").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("U2F5YSBtb2hvbiBtYWFmLCBzYWF0IGluaSBzYXlhIHRpZGFrIGRhcGF0IG1lbWJlcmlrYW4gZGFmdGFyIHNla29sYWggeWFuZyBiZXJwYXJ0aXNpcGFzaSBkYWxhbSBhcGxpa2FzaSAia2FtYmluZ2t1IiBrYXJlbmEga2V0ZXJiYXRhc2FuIGFrc2VzIHRlcmhhZGFwIGRhdGFzZXQgeWFuZyByZWxldmFuLiBEYXRhIHlhbmcgZGlwZXJsdWthbiB0aWRhayBkYXBhdCBkaXRlbXVrYW4gZGFsYW0gcGVuY2FyaWFuIHNlYmVsdW1ueWEsIGRhbiB0YW1wYWtueWEgYWRhIG1hc2FsYWggZGVuZ2FuIGFrc2VzIGtlIHN1bWJlciBkYXRhIHlhbmcgdGVwYXQgYXRhdSBrdWVyaSBwZW5jYXJpYW4uIFVudHVrIGluZm9ybWFzaSBsZWJpaCBsYW5qdXQgYXRhdSBwZW1iYXJ1YW4gdGVya2luaSwgc2F5YSBzYXJhbmthbiB1bnR1ayBtZW5ndW5qdW5naSBzaXR1cyB3ZWIgcmVzbWkgYXBsaWthc2kgImthbWJpbmdrdSIgYXRhdSBtZW5naHVidW5naSBwZW5nZW1iYW5nIGFwbGlrYXNpIHRlcnNlYnV0LiBKaWthIGFkYSBzdW1iZXIgZGF5YSBhdGF1IGluZm9ybWFzaSB0YW1iYWhhbiB5YW5nIGRhcGF0IGRpYWtzZXMsIHNheWEgYWthbiBkZW5nYW4gc2VuYW5nIGhhdGkgbWVtYmFudHUgbGViaWggbGFuanV0Lg==")
$ws.Range("C6").Value = [System.Text.Encoding]::UTF8.GetString($bytes)

$bytes = [Convert]::FromBase64String("MjAyNS0wOC0wMyAyMToyODozMw==")
$ws.Range("A7").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("YmVyYXBhIGp1bWFsaCBzZWtvbGFoIHlhbmcgdGVyZGFmdGFy")
$ws.Range("B7").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("QmVyaWt1dCBhZGFsYWgganVtbGFoIHNla29sYWggeWFuZyB0ZXJkYWZ0YXIgZGFsYW0gYXBsaWthc2kgImthbWJpbmdrdSIgYmVyZGFzYXJrYW4gZGF0YSB5YW5nIHRlcnNlZGlhOgoxLiAqKlBlcnRhbmlhbiBLZWx1YmlyIFRlcnBhZHUgKFBUSykgUFQgS2hhdHVsaXN0aXdhIE51c2FudGFyYSoqCiAgIC0gTG9rYXNpOiBLZWx1YmlyLCBLZWMuIFRqLiBQYWxhcyBVdGFyYSwgS2FidXBhdGVuIEJ1bHVuZ2FuLCBLYWxpbWFudGFuIFV0YXJhCiAgIC0gS29udGFrOiAwODUzNDg3ODEyNjUKICAgLSBXZWJzaXRlOiBbUFRLIEdvYXRdKGh0dHBzOi8vcHRrLWdvYXQuZ2Flbml0cmFpbmluZy5vcmcpCjIuICoqU01LTiAyIElORFJBTUFZVSoqCiAgIC0gTG9rYXNpOiBKbC4gUmF5YSBQYWJlYW4gVWRpayBOby4xNSwgUGFiZWFudWRpaywgS2VjLiBJbmRyYW1heXUsIEthYnVwYXRlbiBJbmRyYW1heXUsIEphd2EgQmFyYXQgNDUyMTkKICAgLSBLb250YWs6IDA4OTUzODAyNjY2ODIKICAgLSBXZWJzaXRlOiBbU01LTiAyIEluZHJhbWF5dV0oaHR0cHM6Ly9zbWtuMmluZHJhbWF5dS5zY2guaWQvKQozLiAqKkZBUk0gUlVNQUggQkVMQUpBUiBQUkFOQVNJU1dBKioKICAgLSBXZWJzaXRlOiBbR2FlbmldKGh0dHBzOi8vd3d3LmdhZW5pLm9yZykKVG90YWwganVtbGFoIHNla29sYWggeWFuZyB0ZXJkYWZ0YXIgYWRhbGFoICoqdGlnYSBzZWtvbGFoKiouIEluZm9ybWFzaSBpbmkgZGlhbWJpbCBkYXJpIGZpbGUgc2Nob29sLmNzdiB5YW5nIGJlcmlzaSBkYXRhIHNla29sYWgtc2Vrb2xhaCB5YW5nIHRlcmRhZnRhciBkYWxhbSBhcGxpa2FzaSB0ZXJzZWJ1dC4KUmVmZXJlbnNpOgotIFtQVEsgR29hdF0oaHR0cHM6Ly9wdGstZ29hdC5nYWVuaXRyYWluaW5nLm9yZykKLSBbU01LTiAyIEluZHJhbWF5dV0oaHR0cHM6Ly9zbWtuMmluZHJhbWF5dS5zY2guaWQvKQotIFtHYWVuaV0oaHR0cHM6Ly93d3cuZ2Flbmkub3JnKQ==")
$ws.Range("C7").Value = [System.Text.Encoding]::UTF8.GetString($bytes)

$bytes = [Convert]::FromBase64String("MjAyNS0wOC0wMyAyMTozMDoxNw==")
$ws.Range("A8").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("YXBhIGFqYSBqZW5pcyBrYW1iaW5nIHlhbmcgdGVyZGF0YQ==")
$ws.Range("B8").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("IyMjIEplbmlzLUplbmlzIEthbWJpbmcKS2FtYmluZyBhZGFsYWggaGV3YW4gdGVybmFrIHlhbmcgbWVtaWxpa2kgYmVyYmFnYWkgamVuaXMgeWFuZyB0ZXJzZWJhciBkaSBzZWx1cnVoIGR1bmlhLiBCZXJpa3V0IGFkYWxhaCBiZWJlcmFwYSBqZW5pcyBrYW1iaW5nIHlhbmcgdW11bSBkaWtlbmFsOgoxLiAqKkthbWJpbmcgQm9lcioqOiBLYW1iaW5nIGluaSBiZXJhc2FsIGRhcmkgQWZyaWthIFNlbGF0YW4gZGFuIGRpa2VuYWwga2FyZW5hIHBlcnR1bWJ1aGFuIHlhbmcgY2VwYXQgZGFuIHByb2R1a3NpIGRhZ2luZyB5YW5nIHRpbmdnaS4gS2FtYmluZyBCb2VyIG1lbWlsaWtpIHR1YnVoIHlhbmcgYmVzYXIgZGFuIG90b3QgeWFuZyBrdWF0LgoyLiAqKkthbWJpbmcgU2FhbmVuKio6IEJlcmFzYWwgZGFyaSBTd2lzcywga2FtYmluZyBpbmkgdGVya2VuYWwgc2ViYWdhaSBwZW5naGFzaWwgc3VzdSB5YW5nIGJhaWsuIEthbWJpbmcgU2FhbmVuIG1lbWlsaWtpIGJ1bHUgcHV0aWggZGFuIHVrdXJhbiB0dWJ1aCB5YW5nIGJlc2FyLgozLiAqKkthbWJpbmcgRXRhd2EgKEphbW5hcGFyaSkqKjogS2FtYmluZyBpbmkgYmVyYXNhbCBkYXJpIEluZGlhIGRhbiBkaWtlbmFsIGthcmVuYSBrZW1hbXB1YW4gcHJvZHVrc2kgc3VzdSB5YW5nIHRpbmdnaS4gS2FtYmluZyBFdGF3YSBqdWdhIHNlcmluZyBkaWd1bmFrYW4gZGFsYW0gcGVyc2lsYW5nYW4gdW50dWsgbWVuaW5na2F0a2FuIGt1YWxpdGFzIGthbWJpbmcgbG9rYWwuCjQuICoqS2FtYmluZyBLYWNhbmcqKjogSmVuaXMga2FtYmluZyBsb2thbCBJbmRvbmVzaWEgeWFuZyBtZW1pbGlraSB1a3VyYW4gdHVidWgga2VjaWwgZGFuIHRhaGFuIHRlcmhhZGFwIGtvbmRpc2kgbGluZ2t1bmdhbiB5YW5nIGtlcmFzLiBLYW1iaW5nIGluaSBzZXJpbmcgZGlwZWxpaGFyYSB1bnR1ayBkaWFtYmlsIGRhZ2luZ255YS4KNS4gKipLYW1iaW5nIEFuZ2xvLU51YmlhbioqOiBLYW1iaW5nIGluaSBhZGFsYWggaGFzaWwgcGVyc2lsYW5nYW4gYW50YXJhIGthbWJpbmcgSW5nZ3JpcyBkYW4ga2FtYmluZyBkYXJpIFRpbXVyIFRlbmdhaC4gTWVyZWthIGRpa2VuYWwga2FyZW5hIHByb2R1a3NpIHN1c3UgeWFuZyB0aW5nZ2kgZGFuIGt1YWxpdGFzIGRhZ2luZyB5YW5nIGJhaWsuCjYuICoqS2FtYmluZyBBbHBpbmUqKjogQmVyYXNhbCBkYXJpIFBlZ3VudW5nYW4gQWxwZW4sIGthbWJpbmcgaW5pIGRpa2VuYWwga2FyZW5hIHByb2R1a3NpIHN1c3UgeWFuZyBiYWlrIGRhbiBrZW1hbXB1YW4gYmVyYWRhcHRhc2kgZGVuZ2FuIGJlcmJhZ2FpIGtvbmRpc2kgaWtsaW0uCjcuICoqS2FtYmluZyBUb2dnZW5idXJnKio6IEthbWJpbmcgaW5pIGJlcmFzYWwgZGFyaSBTd2lzcyBkYW4gbWVydXBha2FuIHNhbGFoIHNhdHUgcmFzIGthbWJpbmcgcGVyYWggdGVydHVhLiBNZXJla2EgbWVtaWxpa2kgYnVsdSBjb2tsYXQgZGVuZ2FuIHRhbmRhIHB1dGloIGRpIHdhamFoIGRhbiBrYWtpLgpVbnR1ayBpbmZvcm1hc2kgbGViaWggbGFuanV0IG1lbmdlbmFpIGplbmlzLWplbmlzIGthbWJpbmcsIEFuZGEgZGFwYXQgbWVuZ3VuanVuZ2kgW1dpa2lwZWRpYSAtIEdvYXRdKGh0dHBzOi8vZW4ud2lraXBlZGlhLm9yZy93aWtpL0dvYXQpIGF0YXUgc3VtYmVyIGxhaW4geWFuZyB0ZXJwZXJjYXlhLgpSZWZlcmVuc2k6Ci0gW1dpa2lwZWRpYSAtIEdvYXRdKGh0dHBzOi8vZW4ud2lraXBlZGlhLm9yZy93aWtpL0dvYXQp")
$ws.Range("C8").Value = [System.Text.Encoding]::UTF8.GetString($bytes)

$bytes = [Convert]::FromBase64String("MjAyNS0wOC0wMyAyMTozMTo1MQ==")
$ws.Range("A9").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("YmVyaWthbiBha3UgZGF0YSBrYW1iaW5nIGRlbmdhbiBqZW5pcyBldGF3YQ==")
$ws.Range("B9").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("IyMjIEthbWJpbmcgRXRhd2EKS2FtYmluZyBFdGF3YSwganVnYSBkaWtlbmFsIHNlYmFnYWkga2FtYmluZyBKYW1uYXBhcmksIGFkYWxhaCBzYWxhaCBzYXR1IGplbmlzIGthbWJpbmcgeWFuZyBiZXJhc2FsIGRhcmkgSW5kaWEgZGFuIHRlbGFoIGJhbnlhayBkaWJ1ZGlkYXlha2FuIGRpIEluZG9uZXNpYS4gS2FtYmluZyBpbmkgdGVya2VuYWwga2FyZW5hIHVrdXJhbiB0dWJ1aG55YSB5YW5nIGJlc2FyIGRhbiBrZW1hbXB1YW4gcHJvZHVrc2lueWEgeWFuZyBiYWlrLCBiYWlrIHVudHVrIHN1c3UgbWF1cHVuIGRhZ2luZy4KIyMjIyBDaXJpLWNpcmkgS2FtYmluZyBFdGF3YToKLSAqKlVrdXJhbiBUdWJ1aDoqKiBLYW1iaW5nIEV0YXdhIG1lbWlsaWtpIHR1YnVoIHlhbmcgYmVzYXIgZGVuZ2FuIHRpbmdnaSBtZW5jYXBhaSA5MC0xMDAgY20gdW50dWsgamFudGFuIGRhbiA3MC05MCBjbSB1bnR1ayBiZXRpbmEuCi0gKipCZXJhdCBCYWRhbjoqKiBCZXJhdCBrYW1iaW5nIGphbnRhbiBkZXdhc2EgYmlzYSBtZW5jYXBhaSA5MSBrZywgc2VkYW5na2FuIGJldGluYSBzZWtpdGFyIDYzIGtnLgotICoqQnVsdToqKiBNZW1pbGlraSBidWx1IHlhbmcgcGFuamFuZyBkYW4gbGViYXQsIHRlcnV0YW1hIGRpIGJhZ2lhbiBwYWhhIGRhbiBwdW5nZ3VuZy4gV2FybmEgYnVsdSBiZXJ2YXJpYXNpLCB0ZXRhcGkgdW11bW55YSBiZXJ3YXJuYSBwdXRpaCBkZW5nYW4gYmVyY2FrIGNva2xhdCBhdGF1IGhpdGFtLgotICoqVGVsaW5nYToqKiBUZWxpbmdhIHBhbmphbmcgZGFuIG1lbmdnYW50dW5nLCBtZW5qYWRpIHNhbGFoIHNhdHUgY2lyaSBraGFzIGthbWJpbmcgaW5pLgotICoqUHJvZHVrc2kgU3VzdToqKiBLYW1iaW5nIEV0YXdhIGRpa2VuYWwgc2ViYWdhaSBwZW5naGFzaWwgc3VzdSB5YW5nIGJhaWssIGRlbmdhbiBwcm9kdWtzaSBzdXN1IGhhcmlhbiBiaXNhIG1lbmNhcGFpIDMtNCBsaXRlci4KLSAqKkFkYXB0YXNpOioqIEthbWJpbmcgaW5pIG1hbXB1IGJlcmFkYXB0YXNpIGRlbmdhbiBiYWlrIGRpIGJlcmJhZ2FpIGtvbmRpc2kgaWtsaW0sIG1lbWJ1YXRueWEgcG9wdWxlciBkaSBiZXJiYWdhaSBkYWVyYWggZGkgSW5kb25lc2lhLgojIyMjIE1hbmZhYXQ6Ci0gKipTdXN1OioqIFN1c3Uga2FtYmluZyBFdGF3YSBrYXlhIGFrYW4gbnV0cmlzaSBkYW4gc2VyaW5nIGRpZ3VuYWthbiBzZWJhZ2FpIGFsdGVybmF0aWYgc3VzdSBzYXBpLCB0ZXJ1dGFtYSBiYWdpIG1lcmVrYSB5YW5nIGFsZXJnaSB0ZXJoYWRhcCBsYWt0b3NhLgotICoqRGFnaW5nOioqIFNlbGFpbiBzdXN1LCBrYW1iaW5nIGluaSBqdWdhIGRpcGVsaWhhcmEgdW50dWsgZGFnaW5nbnlhIHlhbmcgYmVya3VhbGl0YXMgdGluZ2dpLgotICoqUGVtYmlha2FuOioqIEthbWJpbmcgRXRhd2Egc2VyaW5nIGRpZ3VuYWthbiBkYWxhbSBwcm9ncmFtIHBlbWJpYWthbiB1bnR1ayBtZW5pbmdrYXRrYW4ga3VhbGl0YXMgZ2VuZXRpayBrYW1iaW5nIGxva2FsLgojIyMjIFBlbWVsaWhhcmFhbjoKLSAqKlBha2FuOioqIEthbWJpbmcgRXRhd2EgbWVtZXJsdWthbiBwYWthbiBiZXJrdWFsaXRhcyB0aW5nZ2kgeWFuZyBrYXlhIGFrYW4gcHJvdGVpbiBkYW4gc2VyYXQgdW50dWsgbWVuZHVrdW5nIHBlcnR1bWJ1aGFubnlhLgotICoqS2FuZGFuZzoqKiBLYW5kYW5nIGhhcnVzIGJlcnNpaCBkYW4gY3VrdXAgbHVhcyB1bnR1ayBtZW5kdWt1bmcgYWt0aXZpdGFzIGthbWJpbmcsIHNlcnRhIGRpbGVuZ2thcGkgZGVuZ2FuIHNpc3RlbSBkcmFpbmFzZSB5YW5nIGJhaWsgdW50dWsgbWVuamFnYSBrZWJlcnNpaGFuLgotICoqS2VzZWhhdGFuOioqIFBlcmx1IGRpbGFrdWthbiBwZW1lcmlrc2FhbiBrZXNlaGF0YW4gc2VjYXJhIHJ1dGluIHVudHVrIG1lbmNlZ2FoIHBlbnlha2l0IGRhbiBtZW1hc3Rpa2FuIHByb2R1a3Rpdml0YXMga2FtYmluZyB0ZXRhcCBvcHRpbWFsLgpVbnR1ayBpbmZvcm1hc2kgbGViaWggbGFuanV0IHRlbnRhbmcga2FtYmluZyBFdGF3YSwgQW5kYSBiaXNhIG1lbmd1bmp1bmdpIFtzdW1iZXIgcmVmZXJlbnNpXShodHRwczovL2lkLndpa2lwZWRpYS5vcmcvd2lraS9LYW1iaW5nX0V0YXdhKS4KUmVmZXJlbnNpOgotIFtXaWtpcGVkaWE6IEthbWJpbmcgRXRhd2FdKGh0dHBzOi8vaWQud2lraXBlZGlhLm9yZy93aWtpL0thbWJpbmdfRXRhd2Ep")
$ws.Range("C9").Value = [System.Text.Encoding]::UTF8.GetString($bytes)

$bytes = [Convert]::FromBase64String("MjAyNS0wOC0wNSAwODo1ODo1Nw==")
$ws.Range("A10").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("QXBhIGl0dSBLYW1iaW5n")
$ws.Range("B10").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("IyMjIEFwYSBJdHUgS2FtYmluZz8KS2FtYmluZyBhZGFsYWggaGV3YW4gdGVybmFrIHlhbmcgdGVybWFzdWsgZGFsYW0ga2VsdWFyZ2EgQm92aWRhZSBkYW4gc3ViZmFtaWxpIENhcHJpbmFlLiBOYW1hIGlsbWlhaCBrYW1iaW5nIGFkYWxhaCAqQ2FwcmEgYWVnYWdydXMgaGlyY3VzKi4gS2FtYmluZyB0ZWxhaCBkaWRvbWVzdGlrYXNpIHNlamFrIHJpYnVhbiB0YWh1biB5YW5nIGxhbHUgZGFuIG1lbWlsaWtpIHBlcmFuIHBlbnRpbmcgZGFsYW0ga2VoaWR1cGFuIG1hbnVzaWEsIGJhaWsgc2ViYWdhaSBzdW1iZXIgZGFnaW5nLCBzdXN1LCBrdWxpdCwgbWF1cHVuIGJ1bHUuIApLYW1iaW5nIGRpa2VuYWwgc2ViYWdhaSBoZXdhbiB5YW5nIGFkYXB0aWYgZGFuIGRhcGF0IGhpZHVwIGRpIGJlcmJhZ2FpIGxpbmdrdW5nYW4sIGRhcmkgcGFkYW5nIHJ1bXB1dCBoaW5nZ2EgZGFlcmFoIGJlcmJhdHUuIE1lcmVrYSBtZW1pbGlraSBrZW1hbXB1YW4gbWVydW1wdXQgeWFuZyBiYWlrIGRhbiBzZXJpbmcgZGlndW5ha2FuIHVudHVrIG1lbWJhbnR1IG1lbmdlbmRhbGlrYW4gcGVydHVtYnVoYW4gdmVnZXRhc2kgeWFuZyBiZXJsZWJpaGFuLgpEaSBJbmRvbmVzaWEsIHRlcmRhcGF0IGJlcmJhZ2FpIGplbmlzIGthbWJpbmcgeWFuZyBkaWJ1ZGlkYXlha2FuLCBhbnRhcmEgbGFpbiBrYW1iaW5nIGthY2FuZywga2FtYmluZyBldGF3YSwgZGFuIGthbWJpbmcgcGVyYW5ha2FuIGV0YXdhLiBTZXRpYXAgamVuaXMgbWVtaWxpa2kga2FyYWt0ZXJpc3RpayBkYW4ga2V1bmdndWxhbiB0ZXJzZW5kaXJpLCBzZXBlcnRpIHByb2R1a3NpIHN1c3UgeWFuZyB0aW5nZ2kgYXRhdSBrZW1hbXB1YW4gYWRhcHRhc2kgeWFuZyBiYWlrIHRlcmhhZGFwIGxpbmdrdW5nYW4gbG9rYWwuClVudHVrIGluZm9ybWFzaSBsZWJpaCBsYW5qdXQgbWVuZ2VuYWkga2FtYmluZywgQW5kYSBkYXBhdCBtZW5ndW5qdW5naSBzdW1iZXItc3VtYmVyIGJlcmlrdXQ6Ci0gW0thbWJpbmcgLSBXaWtpcGVkaWEgYmFoYXNhIEluZG9uZXNpYSwgZW5zaWtsb3BlZGlhIGJlYmFzXShodHRwczovL2lkLndpa2lwZWRpYS5vcmcvd2lraS9LYW1iaW5nKQotIFtOYW1hIElsbWlhaCBLYW1iaW5nIGRhbiBKZW5pcy1KZW5pcyBLYW1iaW5nIHlhbmcgQWRhIGRpIEluZG9uZXNpYSAtIEtvbXBhcy5jb21dKGh0dHBzOi8vd3d3LmtvbXBhcy5jb20vc2FpbnMvcmVhZC8yMDIyLzAxLzA3LzE3NDYwMDkyMy9uYW1hLWlsbWlhaC1rYW1iaW5nLWRhbi1qZW5pcy1qZW5pcy1rYW1iaW5nLXlhbmctYWRhLWRpLWluZG9uZXNpYSkKLSBbUEVNVURBIFBFVEVSTkFLIEtBTUJJTkcgQlVOR0tBUkFOIC0gUGFnZXJzYXJpXShodHRwOi8vcGFnZXJzYXJpLXBhdGVhbi5kZXNhLmlkL3BvdGVuc2lkZXRhaWwvZVVWdllsWkZZbWRoWkRVMmREaEpNRUZyU200MGR6MDkvcGVtdWRhLXBldGVybmFrLWthbWJpbmctYnVuZ2thcmFuLmh0bWwp")
$ws.Range("C10").Value = [System.Text.Encoding]::UTF8.GetString($bytes)

$bytes = [Convert]::FromBase64String("MjAyNS0wOC0wNSAwODo1OToyOA==")
$ws.Range("A11").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("QXBhIHlhbmcgY29jb2sgc2ViYWdhaSB0YW1iYWhhbiBwYWthbiB1bnR1ayBrYW1iaW5n")
$ws.Range("B11").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("IyMjIFBha2FuIFRhbWJhaGFuIHlhbmcgQ29jb2sgdW50dWsgS2FtYmluZwpLYW1iaW5nIG1lbWJ1dHVoa2FuIHBha2FuIHlhbmcgc2VpbWJhbmcgdW50dWsgbWVuZHVrdW5nIHBlcnR1bWJ1aGFuIGRhbiBrZXNlaGF0YW5ueWEuIFNlbGFpbiBwYWthbiB1dGFtYSBiZXJ1cGEgcnVtcHV0IGRhbiBoaWphdWFuLCBhZGEgYmViZXJhcGEgamVuaXMgcGFrYW4gdGFtYmFoYW4geWFuZyBkYXBhdCBkaWJlcmlrYW4ga2VwYWRhIGthbWJpbmcgdW50dWsgbWVtZW51aGkga2VidXR1aGFuIG51dHJpc2lueWE6CjEuICoqS29uc2VudHJhdCoqOiBLb25zZW50cmF0IG1lcnVwYWthbiBwYWthbiB0YW1iYWhhbiB5YW5nIHBlbnRpbmcgdW50dWsga2FtYmluZywgdGVydXRhbWEgdW50dWsgbWVuaW5na2F0a2FuIHByb2R1a3NpIHN1c3UgZGFuIGRhZ2luZy4gS29uc2VudHJhdCBzZWJhaWtueWEgZGliZXJpa2FuIGRhbGFtIGJlbnR1ayBrYXNhciBhZ2FyIG11ZGFoIGRpY2VybmEgb2xlaCBrYW1iaW5nLiBbQmFjYSBzZWxlbmdrYXBueWEgZGkgc2luaV0oaHR0cHM6Ly9iYnBwa3VwYW5nLmJwcHNkbXAucGVydGFuaWFuLmdvLmlkL2Jsb2cvdHJpay1tZW1iZXJpa2FuLXBha2FuLWthbWJpbmctcG90b25nLXlhbmctYmVuYXIpLgoyLiAqKkhpamF1YW4gZGFuIERlZGF1bmFuKio6IFBha2FuIGhpamF1YW4gc2VwZXJ0aSBydW1wdXQgZ2FqYWgsIGRhdW4gc2luZ2tvbmcsIGRhbiBrYWxpYW5kcmEgc2FuZ2F0IGJhZ3VzIHVudHVrIGthbWJpbmcga2FyZW5hIGtheWEgYWthbiBzZXJhdCBkYW4gbnV0cmlzaS4gRGVkYXVuYW4gZGFyaSB0YW5hbWFuIGxlZ3VtaW5vc2Egc2VwZXJ0aSB0dXJpIGRhbiBnYW1hbCBqdWdhIGRhcGF0IG1lbmphZGkgcGlsaWhhbiB5YW5nIGJhaWsuIFtCYWNhIHNlbGVuZ2thcG55YSBkaSBzaW5pXShodHRwczovL2ZwcC51bWtvLmFjLmlkLzIwMjEvMDYvMjMvamVuaXMtcGFrYW4taGlqYXVhbi15YW5nLWJhZ3VzLXVudHVrLWthbWJpbmcvKS4KMy4gKipCaWppLWJpamlhbioqOiBCaWppLWJpamlhbiBzZXBlcnRpIGphZ3VuZyBkYW4ga2VkZWxhaSBkYXBhdCBkaWJlcmlrYW4gc2ViYWdhaSBzdW1iZXIgZW5lcmdpIHRhbWJhaGFuLCB0ZXJ1dGFtYSB1bnR1ayBhbmFrYW4ga2FtYmluZyB5YW5nIG1lbWJ1dHVoa2FuIGxlYmloIGJhbnlhayBlbmVyZ2kgdW50dWsgcGVydHVtYnVoYW4uIFtCYWNhIHNlbGVuZ2thcG55YSBkaSBzaW5pXShodHRwczovL3d3dy5hZ3JvcHVzdGFrYS5pZC9rYWJhci90aXBzLXBlbWJlcmlhbi1wYWthbi1rYW1iaW5nLWRhbi1kb21iYS8pLgo0LiAqKkFtcGFzIEhhc2lsIFBlcnRhbmlhbioqOiBBbXBhcyBzZXBlcnRpIGFtcGFzIHRhaHUgZGFuIGFtcGFzIGJpciBkYXBhdCBkaWd1bmFrYW4gc2ViYWdhaSBwYWthbiB0YW1iYWhhbiB5YW5nIGVrb25vbWlzIGRhbiBrYXlhIG51dHJpc2kuIFtCYWNhIHNlbGVuZ2thcG55YSBkaSBzaW5pXShodHRwczovL2dkbS5pZC9wYWthbi1rYW1iaW5nLXNlbGFpbi1ydW1wdXQvKS4KNS4gKipGZXJtZW50YXNpIFBha2FuKio6IEZlcm1lbnRhc2kgcGFrYW4gZGFwYXQgbWVuaW5na2F0a2FuIGtlY2VybmFhbiBkYW4gbmlsYWkgZ2l6aSBwYWthbi4gSW5pIGJpc2EgbWVuamFkaSBwaWxpaGFuIGJhaWsgc2FhdCBoaWphdWFuIHN1bGl0IGRpZGFwYXRrYW4uIFtCYWNhIHNlbGVuZ2thcG55YSBkaSBzaW5pXShodHRwczovL2dkbS5pZC9wYWthbi1rYW1iaW5nLXNlbGFpbi1ydW1wdXQvKS4KRGVuZ2FuIG1lbWJlcmlrYW4gcGFrYW4gdGFtYmFoYW4geWFuZyB0ZXBhdCwgcGV0ZXJuYWsgZGFwYXQgbWVtYXN0aWthbiBrYW1iaW5nIG1lbmRhcGF0a2FuIG51dHJpc2kgeWFuZyBjdWt1cCB1bnR1ayBwZXJ0dW1idWhhbiwgcHJvZHVrc2ksIGRhbiBrZXNlaGF0YW5ueWEuIFBhc3Rpa2FuIHVudHVrIHNlbGFsdSBtZW55ZXN1YWlrYW4gamVuaXMgZGFuIGp1bWxhaCBwYWthbiBkZW5nYW4ga2VidXR1aGFuIHNwZXNpZmlrIGthbWJpbmcgQW5kYS4=")
$ws.Range("C11").Value = [System.Text.Encoding]::UTF8.GetString($bytes)

$bytes = [Convert]::FromBase64String("MjAyNS0wOC0wNSAwOTowMTo0Mg==")
$ws.Range("A12").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("a2FtYmluZyBqZW5pcyBhcGEgeWFuZyBwYWxpbmcgYmFueWFrIGRpIGthbWJpbmdrdT8=")
$ws.Range("B12").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("QmVyZGFzYXJrYW4gaW5mb3JtYXNpIHlhbmcgdGVyc2VkaWEsIHNheWEgYmVsdW0gZGFwYXQgbWVuZW11a2FuIGRhdGEgc3Blc2lmaWsgbWVuZ2VuYWkgamVuaXMga2FtYmluZyB5YW5nIHBhbGluZyBiYW55YWsgZGkgYXBsaWthc2kgIkthbWJpbmdrdSIuIERhdGEgeWFuZyByZWxldmFuIG11bmdraW4gdGVyc2ViYXIgZGkgYmVyYmFnYWkgZmlsZSBzZXBlcnRpIGBnb2F0LmNzdmAgZGFuIGB0eXBlLmNzdmAuIFVudHVrIG1lbmRhcGF0a2FuIGphd2FiYW4geWFuZyBsZWJpaCBha3VyYXQsIGRpc2FyYW5rYW4gdW50dWsgbWVsYWt1a2FuIGFuYWxpc2lzIGxhbmdzdW5nIHBhZGEgZGF0YSB5YW5nIHRlcmRhcGF0IGRhbGFtIGZpbGUtZmlsZSB0ZXJzZWJ1dCwgZGVuZ2FuIGZva3VzIHBhZGEga29sb20geWFuZyBtZW51bmp1a2thbiBqZW5pcyBrYW1iaW5nIGRhbiBqdW1sYWhueWEuIEFuYWxpc2lzIGluaSBkYXBhdCBkaWxha3VrYW4gZGVuZ2FuIG1lbmdndW5ha2FuIHBlcmFuZ2thdCBsdW5hayBwZW5nb2xhaCBkYXRhIHNlcGVydGkgTWljcm9zb2Z0IEV4Y2VsIGF0YXUgR29vZ2xlIFNoZWV0cyB1bnR1ayBtZW5naWRlbnRpZmlrYXNpIGplbmlzIGthbWJpbmcgeWFuZyBwYWxpbmcgYmFueWFrIGJlcmRhc2Fya2FuIGp1bWxhaCBrZW11bmN1bGFubnlhIGRhbGFtIGRhdGFzZXQgdGVyc2VidXQuClNlYmFnYWkgbGFuZ2thaCBhd2FsLCBBbmRhIGRhcGF0IG1lbWJ1a2EgZmlsZSBgZ29hdC5jc3ZgIGRhbiBgdHlwZS5jc3ZgLCBrZW11ZGlhbiBtZW5jYXJpIGtvbG9tIHlhbmcgcmVsZXZhbiBzZXBlcnRpICJqZW5pcyBrYW1iaW5nIiBkYW4gImp1bWxhaCIuIFNldGVsYWggaXR1LCBsYWt1a2FuIHBlbmdlbG9tcG9rYW4gZGFuIHBlbmdoaXR1bmdhbiB1bnR1ayBtZW5lbXVrYW4gamVuaXMga2FtYmluZyB5YW5nIHBhbGluZyBiYW55YWsuIEppa2EgQW5kYSBtZW1lcmx1a2FuIHBhbmR1YW4gbGViaWggbGFuanV0IHRlbnRhbmcgY2FyYSBtZWxha3VrYW4gYW5hbGlzaXMgZGF0YSBpbmksIEFuZGEgZGFwYXQgbWVuZ3VuanVuZ2kgc3VtYmVyIGRheWEgb25saW5lIGF0YXUgdHV0b3JpYWwgeWFuZyB0ZXJzZWRpYS4KVW50dWsgaW5mb3JtYXNpIGxlYmloIGxhbmp1dCBkYW4gcGFuZHVhbiB0ZW50YW5nIGNhcmEgbWVsYWt1a2FuIGFuYWxpc2lzIGRhdGEsIEFuZGEgZGFwYXQgbWVuZ3VuanVuZ2kgW3R1dG9yaWFsIGFuYWxpc2lzIGRhdGEgZGVuZ2FuIEV4Y2VsXShodHRwczovL3N1cHBvcnQubWljcm9zb2Z0LmNvbS9pZC1pZC9leGNlbCkgYXRhdSBbR29vZ2xlIFNoZWV0c10oaHR0cHM6Ly9zdXBwb3J0Lmdvb2dsZS5jb20vZG9jcy9hbnN3ZXIvMzA5MzQ4MD9obD1pZCkuClJlZmVyZW5zaToKLSBbTWljcm9zb2Z0IEV4Y2VsIFN1cHBvcnRdKGh0dHBzOi8vc3VwcG9ydC5taWNyb3NvZnQuY29tL2lkLWlkL2V4Y2VsKQotIFtHb29nbGUgU2hlZXRzIFN1cHBvcnRdKGh0dHBzOi8vc3VwcG9ydC5nb29nbGUuY29tL2RvY3MvYW5zd2VyLzMwOTM0ODA/aGw9aWQp")
$ws.Range("C12").Value = [System.Text.Encoding]::UTF8.GetString($bytes)

$bytes = [Convert]::FromBase64String("MjAyNS0wOC0wNSAxNDoxNDo1Nw==")
$ws.Range("A13").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("QmVyaWthbiBha3UgcmVrb21lbmRhc2kgcGFrYW4gc2VsYWluIGhpamF1YW4gdW50dWsga2FtYmluZw==")
$ws.Range("B13").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
$bytes = [Convert]::FromBase64String("QmVyaWt1dCBhZGFsYWggYmViZXJhcGEgcmVrb21lbmRhc2kgcGFrYW4gYWx0ZXJuYXRpZiB1bnR1ayBrYW1iaW5nIHNlbGFpbiBoaWphdWFuOgoxLiAqKkRhaGFuIGRhbiBCYXRhbmcqKjogS2FtYmluZyBtZW55dWthaSBtYWthbiBkYWhhbiBkYW4gYmF0YW5nLCB5YW5nIGRhcGF0IGRpdGVtcGF0a2FuIGRpIHRlbXBhdCBwYWthbiBqZXJhbWkuIEluaSB0ZXJtYXN1ayBidWFoIHN1bWFjIGRhbiBkYWhhbiBwb2hvbiBsYWlubnlhLgogICAtIFN1bWJlcjogW1RoZSBSZW5haXNzYW5jZSBIb3VzZXdpZmVdKGh0dHBzOi8vdGhlcmVuYWlzc2FuY2Vob3VzZXdpZmUud2VlYmx5LmNvbS9ob21lL2FsdGVybmF0aXZlLWdvYXQtZmVlZCkKMi4gKipKYWd1bmcqKjogU2VyaW5nIGRpc2VidXQgc2ViYWdhaSAicGVybWVuIGthbWJpbmcsIiBqYWd1bmcgYWRhbGFoIHN1cGxlbWVuIGVuZXJnaSB5YW5nIHNhbmdhdCBiYWlrLiBKYWd1bmcgbXVyYWggZGFuIGRhcGF0IGRpZ3VuYWthbiB1bnR1ayBtZW5ndW1wdWxrYW4sIG1lbWluZGFoa2FuLCBkYW4gbWVuZ2FuZGFuZ2thbiBrYW1iaW5nLgogICAtIFN1bWJlcjogW1ZhbCBWZXJkZSBDb3VudHldKGh0dHBzOi8vdmFsdmVyZGUuYWdyaWxpZmUub3JnL2ZpbGVzLzIwMjAvMTEvU3VwcGxlbWVudGF0aW9uLW9mLUdyYXppbmctR29hdHMucGRmKQozLiAqKkFtcGFzIEJpdCBhdGF1IEt1bGl0IEtlZGVsYWkqKjogSW5pIGFkYWxhaCBzdXBsZW1lbiBzZXJhdCB5YW5nIGRhcGF0IGRpZmVybWVudGFzaSB5YW5nIGRhcGF0IGRpZ3VuYWthbiBrZXRpa2Egc3VtYmVyIGVuZXJnaSB0YW1iYWhhbiBkaWJ1dHVoa2FuLCBzZXBlcnRpIHNlbGFtYSBrZWhhbWlsYW4gYXRhdSBhd2FsIGxha3Rhc2kuCiAgIC0gU3VtYmVyOiBbTWVyY2sgVmV0ZXJpbmFyeSBNYW51YWxdKGh0dHBzOi8vd3d3Lm1lcmNrdmV0bWFudWFsLmNvbS9tYW5hZ2VtZW50LWFuZC1udXRyaXRpb24vcHJldmVudGF0aXZlLWhlYWx0aC1jYXJlLWFuZC1odXNiYW5kcnktb2YtZ29hdHMvbnV0cml0aW9uLW9mLWdvYXRzKQo0LiAqKktvbnNlbnRyYXQqKjogQmlqaS1iaWppYW4gYXRhdSBrb25zZW50cmF0IGhhcnVzIGRpYmVyaWthbiBrZXBhZGEga2FtYmluZyBkZW5nYW4ga2VidXR1aGFuIGVuZXJnaSB5YW5nIGxlYmloIHRpbmdnaSwgc2VwZXJ0aSBhbmFrIGthbWJpbmcgeWFuZyBzZWRhbmcgdHVtYnVoLCBrYW1iaW5nIGhhbWlsLCBhdGF1IG1lbnl1c3VpLiBQZW50aW5nIHVudHVrIG1lbnllaW1iYW5na2FuIGluaSBkZW5nYW4ga2VidXR1aGFuIG51dHJpc2kgc3Blc2lmaWsga2FtYmluZy4KICAgLSBTdW1iZXI6IFtVbml2ZXJzaXR5IG9mIFRlbm5lc3NlZV0oaHR0cHM6Ly92ZXRtZWQudGVubmVzc2VlLmVkdS93cC1jb250ZW50L3VwbG9hZHMvc2l0ZXMvNC9VVENWTV9MQUNTLUZlZWRpbmdHb2F0cy5wZGYpCjUuICoqQmlqaSBLYXBhcyBVdHVoIGRhbiBUZXB1bmcgS2VkZWxhaSoqOiBJbmkgYWRhbGFoIHBpbGloYW4gcGFrYW4gdGluZ2dpIHByb3RlaW4geWFuZyBkYXBhdCBkaWd1bmFrYW4gdW50dWsga2FtYmluZyBiZXRpbmEgbWVueXVzdWkgZGFuIGthbWJpbmcgbGFpbm55YSB5YW5nIG1lbWJ1dHVoa2FuIGFzdXBhbiBwcm90ZWluIGxlYmloIHRpbmdnaS4KICAgLSBTdW1iZXI6IFtOb3J0aCBDYXJvbGluYSBTdGF0ZSBVbml2ZXJzaXR5XShodHRwczovL2NvbnRlbnQuY2VzLm5jc3UuZWR1L251dHJpdGlvbmFsLWZlZWRpbmctbWFuYWdlbWVudC1vZi1tZWF0LWdvYXRzKQpQaWxpaGFuLXBpbGloYW4gaW5pIG1lbnllZGlha2FuIGJlcmJhZ2FpIG51dHJpc2kgZGFuIGRhcGF0IGRpZ3VuYWthbiB1bnR1ayBtZWxlbmdrYXBpIGF0YXUgbWVuZ2dhbnRpa2FuIHBha2FuIGhpamF1YW4gdHJhZGlzaW9uYWwgdGVyZ2FudHVuZyBwYWRhIGtlYnV0dWhhbiBzcGVzaWZpayBrYW1iaW5nLg==")
$ws.Range("C13").Value = [System.Text.Encoding]::UTF8.GetString($bytes)
